$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update is_active (column D) to FALSE for the "Others" gender rows across all languages
$ws.Range("D4").Value = $false
$ws.Range("D7").Value = $false
$ws.Range("D10").Value = $false

# Update the last selected cell in the sheet (matches active cell after edits)
$ws.Range("D12").Select()
